# Update cryptocurrency Price (D) and Volume(1h) (E) columns per the
# latest scrape. Values that look numeric but carry meaningful trailing
# zeros (e.g. 138.30, 5.00) are written with a leading apostrophe so
# Excel keeps them as literal text instead of normalising the digits
# (dropping the trailing zero / switching to scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.027.19'
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '2.377.98'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '560.59'
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").Value = '''138.30'
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").Value = '''0.530'
$ws.Range("E8").Value = '  +0.51%  '
$ws.Range("D9").Value = '2.375.69'
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("D10").Value = '0.106'
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '''5.10'
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D14").Value = '25.74'
$ws.Range("E14").Value = '  -0.85%  '
$ws.Range("D15").Value = '2.794.33'
$ws.Range("D16").Value = '''0.0000166'
$ws.Range("E16").Value = '  -2.67%  '
$ws.Range("D17").Value = '59.859.14'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").Value = '2.376.07'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '8.11'
$ws.Range("E19").Value = '  +13.07%  '
$ws.Range("D20").Value = '10.53'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").Value = '321.41'
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("D22").Value = '4.06'
$ws.Range("E22").Value = '  +1.10%  '
$ws.Range("D23").Value = '6.04'
$ws.Range("E23").Value = '  -1.02%  '
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").Value = '1.82'
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("D26").Value = '64.04'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").Value = '559.07'
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").Value = '8.16'
$ws.Range("E28").Value = '  -6.14%  '
$ws.Range("D30").Value = '0.0₃0928'
$ws.Range("E30").Value = '  +2.18%  '
$ws.Range("E31").Value = '  +3.33%  '
$ws.Range("D32").Value = '1.31'
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("E33").Value = '  -2.63%  '
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("E36").Value = '  +4.24%  '
$ws.Range("D37").Value = '153.42'
$ws.Range("E37").Value = '  +4.19%  '
$ws.Range("D38").Value = '0.367'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("D39").Value = '4.56'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("D40").Value = '18.15'
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").Value = '''5.00'
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = '''41.50'
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").Value = '1.65'
$ws.Range("E44").Value = '  -1.18%  '
$ws.Range("D45").Value = '2.43'
$ws.Range("E45").Value = '  +3.53%  '
$ws.Range("D46").Value = '0.0₆0299'
$ws.Range("E46").Value = '  +4.47%  '
$ws.Range("D47").Value = '140.21'
$ws.Range("E47").Value = '  -0.07%  '
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("D49").Value = '0.585'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '0.0501'
$ws.Range("E50").Value = '  -0.51%  '
$ws.Range("D51").Value = '19.11'
$ws.Range("E51").Value = '  -1.27%  '
